$wb = $excel.ActiveWorkbook

# --- Hoja1: update the "Conversión del día" text cell ---
$ws1 = $wb.Worksheets.Item("Hoja1")
$ws1.Range("A1").Value = "Conversión del día 💰`n✅ Dólar paralelo: 68`n`nBinance`n✅ 1000 Bs = 4.62 = 18299.45 pesos`n✅ 18299.45 pesos = 4.61 = 949.52 Bs`n`nPromedio competencia`n✅ Tasa pesos: 20`n✅ Tasa Bs: 20`n✅ % Ganancia: 20%"

# --- tasas: update N10, O10, N12, O12 ---
$ws2 = $wb.Worksheets.Item("tasas")
$ws2.Range("N10").Value = 216.4
$ws2.Range("O10").Value = 3960
$ws2.Range("N12").Value = 3970.08
$ws2.Range("O12").Value = 206
